$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: add date columns C1, D1, E1 matching B1's look (bold / centered / bordered, literal text) ---
# Assign as a formula returning a literal string first so Excel does not
# auto-convert the "YYYY-MM-DD" text into a date serial number, convert the
# formula result back into a static value, then copy B1's formatting onto it
# so the cell ends up sharing the exact same style as the rest of the header row.
$ws.Cells.Item(1, 3).Formula = '="2021-11-09"'
$ws.Range("C1").Copy()
$ws.Range("C1").PasteSpecial(-4163)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$ws.Cells.Item(1, 4).Formula = '="2021-11-10"'
$ws.Range("D1").Copy()
$ws.Range("D1").PasteSpecial(-4163)
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Cells.Item(1, 5).Formula = '="2021-11-11"'
$ws.Range("E1").Copy()
$ws.Range("E1").PasteSpecial(-4163)
$ws.Range("B1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# --- Column E: "NA" for every data row, except row 13 which holds the numeric value 13 ---
for ($r = 2; $r -le 25; $r++) {
    if ($r -eq 13) {
        $ws.Cells.Item($r, 5).Value = 13
    } else {
        $ws.Cells.Item($r, 5).Value = "NA"
    }
}

# --- Column D numeric exceptions ---
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(19, 4).Value = 13

# --- Column C numeric exception ---
$ws.Cells.Item(21, 3).Value = -1

Write-Host "applied people-counting-summary edit"
